$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.304.70'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '3.594.01'
$ws.Range("E3").Value = '  -2.46%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = "'193.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = "'574.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.83%  '
$ws.Range("D7").Value = '3.588.30'
$ws.Range("E7").Value = '  -2.42%  '
$ws.Range("D8").Value = "'0.617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").Value = "'0.682"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.61%  '
$ws.Range("D11").Value = "'56.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.04%  '
$ws.Range("D12").Value = "'0.150"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.52%  '
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.30%  '
$ws.Range("D14").Value = "'9.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.77%  '
$ws.Range("D15").Value = '4.175.50'
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").Value = '3.590.36'
$ws.Range("E16").Value = '  -2.65%  '
$ws.Range("D17").Value = "'0.126"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").Value = "'18.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.13%  '
$ws.Range("D19").Value = '67.211.77'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").Value = "'12.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.43%  '
$ws.Range("D21").Value = "'1.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.78%  '
$ws.Range("D22").Value = "'402.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("D23").Value = "'4.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.04%  '
$ws.Range("D24").Value = "'86.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.17%  '
$ws.Range("D25").Value = "'11.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("D26").Value = "'2.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.58%  '
$ws.Range("D27").Value = "'12.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").Value = "'6.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("D29").Value = "'3.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.02%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = "'7.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.92%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'8.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.45%  '
$ws.Range("D32").Value = "'31.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.92%  '
$ws.Range("D33").Value = "'639.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.53%  '
$ws.Range("D34").Value = "'12.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.85%  '
$ws.Range("D35").Value = "'0.115"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.28%  '
$ws.Range("D36").Value = "'64.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.47%  '
$ws.Range("D37").Value = "'42.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.22%  '
$ws.Range("D38").Value = "'0.403"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").Value = '0.0₃0773'
$ws.Range("E40").Value = '  -4.07%  '
$ws.Range("D41").Value = '3.232.70'
$ws.Range("E41").Value = '  +13.13%  '
$ws.Range("D42").Value = "'0.134"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Value = "'2.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.26%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").Value = "'0.0417"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.64%  '
$ws.Range("E47").Value = '  +1.85%  '
$ws.Range("E48").Value = '  -5.42%  '
$ws.Range("D49").Value = "'142.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = "'8.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.88%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = "'2.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.63%  '
